# Apply the "up to dual learning step" edit: add eight small vertex-label
# textboxes (V1, V2, V'3, V4, V'5, V'3, V4, V'5) onto the dual-learning
# figure slide, mirroring the pre-existing V1..V5 / V3..V5 label group.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

$EMU_PER_PT = 12700.0

# U+2019 RIGHT SINGLE QUOTATION MARK, used by the author as the "prime" mark.
$rsquo = [char]0x2019

# off-x, off-y, cx, cy (all EMU), text
$labels = @(
    @(7185540, 1635776, 502894, 369332, "V1"),
    @(7331143, 2275640, 502894, 369332, "V2"),
    @(8114288, 2182284, 502894, 369332, ("V" + $rsquo + "3")),
    @(7761420, 2692154, 502894, 369332, "V4"),
    @(8203367, 2872772, 502894, 369332, ("V" + $rsquo + "5")),
    @(8059244, 3211424, 502894, 369332, ("V" + $rsquo + "3")),
    @(7694068, 3747931, 502894, 369332, "V4"),
    @(8214813, 3534331, 502894, 369332, ("V" + $rsquo + "5"))
)

foreach ($lab in $labels) {
    $left   = $lab[0] / $EMU_PER_PT
    $top    = $lab[1] / $EMU_PER_PT
    $width  = $lab[2] / $EMU_PER_PT
    $height = $lab[3] / $EMU_PER_PT
    $text   = $lab[4]

    $tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
    $tb.TextFrame.WordWrap = -1
    $tb.TextFrame.AutoSize = 1
    $tb.Fill.Visible = 0
    $tb.TextFrame.TextRange.Text = $text
}
